$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column header (R2) - "TIME (TOTAL)" is a brand new shared string
$ws.Range("R2").Value = "TIME (TOTAL)"

# New TIME (TOTAL) values for rows 3-17 (kernel 1 en 2 resultaten)
$ws.Range("R3").Value = 114.11
$ws.Range("R4").Value = 129.91
$ws.Range("R5").Value = 122.11
$ws.Range("R6").Value = 144.16
$ws.Range("R7").Value = 115.69
$ws.Range("R8").Value = 224.89
$ws.Range("R9").Value = 164.03
$ws.Range("R10").Value = 145.1
$ws.Range("R11").Value = 123.16
$ws.Range("R12").Value = 1062.86
$ws.Range("R13").Value = 126.04
$ws.Range("R14").Value = 152.89
$ws.Range("R15").Value = 137.46
$ws.Range("R16").Value = 1033.09
$ws.Range("R17").Value = 1083.43

# Auto-size the new column to fit its contents, like Excel does for the
# other bestFit columns on this sheet
$ws.Columns("R").AutoFit()

# Move the active selection to follow the newly added column, matching
# where the cursor ends up after filling in column R
[void]$ws.Range("R18").Select()
